# Improve quality of querymodel testarc
# Rework the "MS" worksheet (injection volume parameter columns) of the
# isa.assay.xlsx workbook:
#   - remove the (now unused) generic "Term Source REF (MS:1001808)" /
#     "Term Accession Number (MS:1001808)" constant values from columns C/D
#   - rename "Parameter [injection volume]" -> "Parameter [injection volume setting]"
#   - give the trailing "Term Source REF ()" / "Term Accession Number ()" columns
#     a concrete ontology reference (AFR:0001577) and fill them in with real
#     term data (source "UO", accession hyperlinked to the UO term)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Columns C (Term Source REF (MS:1001808)) and D (Term Accession Number
# (MS:1001808)) no longer carry a per-row "user-specific" placeholder value.
$ws.Range("C2:D7").ClearContents()

# Column G becomes the concrete term-source for the injection volume setting.
$ws.Range("G2").Value = "UO"
$ws.Range("G3").Value = "UO"
$ws.Range("G4").Value = "UO"
$ws.Range("G5").Value = "UO"
$ws.Range("G6").Value = "UO"
$ws.Range("G7").Value = "UO"

# Rename the injection-volume parameter header.
$ws.Range("E1").Value = "Parameter [injection volume setting]"

# Column H gets the term accession number, expressed as a hyperlink to the
# ontology term.
$ws.Hyperlinks.Add($ws.Range("H2"), "http://purl.obolibrary.org/obo/UO_0000101", "", "", "http://purl.obolibrary.org/obo/UO_0000101")
$ws.Hyperlinks.Add($ws.Range("H3"), "http://purl.obolibrary.org/obo/UO_0000101", "", "", "http://purl.obolibrary.org/obo/UO_0000101")
$ws.Hyperlinks.Add($ws.Range("H4"), "http://purl.obolibrary.org/obo/UO_0000101", "", "", "http://purl.obolibrary.org/obo/UO_0000101")
$ws.Hyperlinks.Add($ws.Range("H5"), "http://purl.obolibrary.org/obo/UO_0000101", "", "", "http://purl.obolibrary.org/obo/UO_0000101")
$ws.Hyperlinks.Add($ws.Range("H6"), "http://purl.obolibrary.org/obo/UO_0000101", "", "", "http://purl.obolibrary.org/obo/UO_0000101")
$ws.Hyperlinks.Add($ws.Range("H7"), "http://purl.obolibrary.org/obo/UO_0000101", "", "", "http://purl.obolibrary.org/obo/UO_0000101")

# Give the two trailing "Term Source REF ()" / "Term Accession Number ()"
# headers their concrete ontology annotation (AFR:0001577).
$ws.Range("G1").Value = "Term Source REF (AFR:0001577)  "
$ws.Range("H1").Value = "Term Accession Number (AFR:0001577)  "

# Move the current selection.
$ws.Range("E10").Select()
